$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing German (informal) INTRO text wording (row 3, column C)
$ws.Range("C3").Value = 'In diesem Test hörst du einige kurze Rhythmen, die du dem dazugehörigen Bild zuordnen sollst. Die Rhythmen werden von jeweils vier Metronomklicks umrahmt. Zunächst hörst du einige Beispiele und machst ein paar Übungsaufgaben.'

# Populate new column D (DE_F - formal German) and column E (IT - Italian)
$ws.Range("D1").Value = 'DE_F'
$ws.Range("E1").Value = 'IT'
$ws.Range("D2").Value = 'Willkommen zum Rhythmuswahrnehmungstest'
$ws.Range("E2").Value = 'Benvenuto/a al test sul ritmo musicale!'
$ws.Range("D3").Value = 'In diesem Test hören Sie einige kurze Rhythmen, die Sie dem dazugehörigen Bild zuordnen sollen. Die Rhythmen werden von jeweils vier Metronomklicks umrahmt. Zunächst hören Sie einige Beispiele und machen ein paar Übungsaufgaben.'
$ws.Range("E3").Value = 'In questo test ascolterai una serie di brevi sequenze ritmiche che dovrai abbinare a un''immagine corrispondente. Ciascuno dei ritmi sarà preceduto e seguito da quattro click di metronomo. Per prima cosa ascolterai alcuni esempi e farai alcune prove per esercitarti.'
$ws.Range("D4").Value = '**Beispiel 1**\\Jeder Rhythmus besteht aus Bassdrum und Claps.\\Die Claps werden als hellblaue Quadrate in der oberen Reihe dargestellt und die Bassdrumschläge als dunkelblaue Quadrate in der unteren Reihe.\\Klicken Sie auf das Abspielzeichen, um den Rhythmus zu hören. Vor und nach dem Rhythmus werden 4 Metronomklicks zu hören sein.'
$ws.Range("E4").Value = '**Esempio 1** Ogni ritmo è composto da suoni di applausi e grancassa \\ Gli applausi sono rappresentati da quadrati azzurri nella riga superiore e la grancassa da quadrati blu scuri nella riga inferiore\\ Ascolta questo breve ritmo facendo click sul pulsante Play. Sentirai 4 click di metronomo prima e dopo il ritmo.'
$ws.Range("D5").Value = '**Clap** **Clap** Bassdrum Bassdrum'
$ws.Range("E5").Value = '**Applauso** **Applauso** Grancassa Grancassa'
$ws.Range("D6").Value = '**Beispiel 2**\\Hier ist ein anderer Rhythmus:'
$ws.Range("E6").Value = '**Esempio 2**\\Ecco qui un altro ritmo:'
$ws.Range("D7").Value = 'Schwach **Stark** Schwach **Stark**'
$ws.Range("E7").Value = 'Debole **Forte** Debole **Forte**'
$ws.Range("D8").Value = 'Rhythmuswahrnehmungs-IQ'
$ws.Range("E8").Value = 'IQ di Percezione Ritmica'
$ws.Range("D9").Value = '**Übungsaufgabe 1**\\Hier ein Beispiel mit vier Elementen.'
$ws.Range("E9").Value = '**Esercitazione  1** \\ Ecco un esempio pratico di quattro suoni'
$ws.Range("D10").Value = 'Welches Bild passt zu dem Rhythmus, den Sie gerade gehört haben? Klicken Sie auf das richtige. Wenn Sie sich nicht sicher sind, dann entscheiden Sie aus dem Bauch heraus.'
$ws.Range("E10").Value = 'Quale immagine corrisponde al ritmo che hai appena sentito? Clicca su quella corretta. Se non lo sai, prova a fare del tuo meglio!'
$ws.Range("D11").Value = 'Erste'
$ws.Range("E11").Value = 'Il primo'
$ws.Range("D12").Value = 'Richtig!'
$ws.Range("E12").Value = 'Giusto!'
$ws.Range("D13").Value = '**Übungsaufgabe 2**\\**{{feedback}}**\\Noch ein Beispielfrage, diesmal mit acht Elementen.'
$ws.Range("E13").Value = '**Esercitazione 2**\\**{{feedback}}**\\ Proviamo con un esempio pratico di 8 suoni'
$ws.Range("D14").Value = 'Falsch.'
$ws.Range("E14").Value = 'Sbagliato.'
$ws.Range("D15").Value = '**{{feedback}}**\\Klicke Sie auf  ''Zurück'', um die Anweisungen erneut zu lesen und die Beispiele erneut zu versuchen,\\ oder klicken Sie auf ''Weiter'', um zum  Haupttest zu gelangen.'
$ws.Range("E15").Value = 'Premi "Indietro" per leggere le istruzioni e ripeti le prove di esercitazione,\\oppure premi "Continua" per procedere con il test principale.'
$ws.Range("D16").Value = 'Zweite'
$ws.Range("E16").Value = 'Il secondo'
$ws.Range("D17").Value = 'Nun geht es mit dem Haupttest los, in dem Ihre  Ergebnisse gespeichert werden.\\Ab jetzt bekommen Sie keine Rückmeldung mehr. Viel Erfolg!'
$ws.Range("E17").Value = 'Stai per iniziare il test principale, in cui i tuoi risultati verranno registrati\\Potrai ascoltare ogni ritmo solo una volta\\ Non riceverai i risultati dopo le singole domande. Buona fortuna!'
$ws.Range("D18").Value = '**Frage {{num_question}} von {{test_length}}**'
$ws.Range("E18").Value = 'Domanda {{num_question}} di {{test_length}}'
$ws.Range("D19").Value = 'Ihr Browser unterstützt kein Audio. Dieser Test funktioniert nicht ohne Audio, sorry!'
$ws.Range("E19").Value = 'Il tuo browser non supporta l’audio. Questo test non può funzionare senz’audio, ci dispiace!'
$ws.Range("D20").Value = 'Nur eins der Bilder passt zu dem Rhythmus. Welches? Wenn Sie nicht sicher sind, dann entscheiden Sie aus dem Bauch heraus.'
$ws.Range("E20").Value = 'Solo una di queste immagini corrisponde al ritmo. Quale? Se non lo sai, prova a fare del tuo meglio!'
$ws.Range("D21").Value = 'Sie haben den Rhythmus-Wahrnehmungs-Test abgeschlossen.\\Von {{num_question}} Aufgaben waren {{num_correct}} richtig.'
$ws.Range("E21").Value = 'Hai completato il test di percezione ritmica!\\Hai risposto correttamente a {{num_correct}} su {{num_question}}.'
$ws.Range("D22").Value = 'Weiter'
$ws.Range("E22").Value = 'Continua'
$ws.Range("D23").Value = 'Welches Bild passt zu dem Rhythmus, den Sie  gerade gehört haben? Klicke Sie auf das richtige. Wenn Sie sich nicht sicher sind, dann entscheiden Sie aus dem Bauch heraus.'
$ws.Range("E23").Value = 'Quale immagine corrisponde al ritmo che hai appena sentito? Clicca su quella corretta. Se non lo sai, prova a fare del tuo meglio!'
$ws.Range("D24").Value = 'Zurück'
$ws.Range("E24").Value = 'Indietro'
$ws.Range("D25").Value = 'Im Vergleich zu anderen, die an dem Test teilgenommen haben, ist Ihr Rhythmus IQ:'
$ws.Range("E25").Value = 'Rispetto alla popolazione generale il tuo QI ritmico è'
$ws.Range("D26").Value = 'Rhythmus IQ'
$ws.Range("E26").Value = 'IQ Ritmico'
$ws.Range("D27").Value = '**Übungsaufgabe 3**\\**{{feedback}}**\\Hier noch eine Aufgabe mit sechszehn Elementen.'
$ws.Range("E27").Value = '**Esercitazione 3**\\**{{feedback}}**\\ Proviamo con un esempio pratico finale di 16 suoni.'
$ws.Range("D28").Value = 'Der Rhythmus wird aus vier, acht oder sechzehn Klängen/Tönen bestehen und Sie hören jeweils vier Metronomschläge vor und nach dem Rhythmus.\\Ihre Aufgabe ist es, den Rhythmus zu hören und dann auf das Bild der vier Bilder klicken, das mit dem Rhyhtmus übereinstimmt, den Sie gerade gehört haben.\\ Lassen Sie uns das mal üben.'
$ws.Range("E28").Value = 'Le sequenze ritmiche saranno composte da quattro, otto o sedici suoni e ci saranno quattro click di metronomo prima e dopo il ritmo effettivo \\ Il tuo compito è ascoltare il ritmo e quindi cliccare sull''immagine tra le quattro opzioni che corrisponde al ritmo che hai appena sentito\\Facciamo un po’ di pratica.'
$ws.Range("D29").Value = 'Klicken Sie hier, falls das Audio nicht spielt.'
$ws.Range("E29").Value = 'Clicca qui se l’audio non parte'
$ws.Range("D30").Value = 'Bitte geben Sie Ihre ID ein'
$ws.Range("E30").Value = 'Per favore inserisci il tuo codice partecipante'
$ws.Range("D31").Value = 'Z.B. 123456'
$ws.Range("E31").Value = 'es. 123456'
$ws.Range("D32").Value = 'Ihre Ergebnisse wurden gespeichert.'
$ws.Range("E32").Value = 'I tuoi risultati sono stati salvati'
$ws.Range("D33").Value = 'Sie können den Browsertab jetzt schließen.'
$ws.Range("E33").Value = 'Adesso puoi chiudere la finestra del browser'
$ws.Range("D34").Value = 'Teilnehmer IDs beginnen mit UK oder AUS, gefolgt von einer einer Zahl, z.B. UK_01 oder AUS_02.'
$ws.Range("E34").Value = 'Il codice partecipante deve iniziare con UK o AUS seguito da un underscore e un numero, es. UK_01 o AUS_02'
$ws.Range("D35").Value = 'Ihr Testergebnis'
$ws.Range("E35").Value = 'Il tuo punteggio'
$ws.Range("D36").Value = 'Werte'
$ws.Range("E36").Value = 'Punteggio'
$ws.Range("D37").Value = 'Rhythmuswahrnehmungstest'
$ws.Range("E37").Value = 'Test sul Ritmo Musicale'
$ws.Range("D38").Value = 'Sie haben den  Rhythmuswahrnehmungstest beendet.'
$ws.Range("E38").Value = 'Hai completato il test sul ritmo musicale.'

# Style header E1: bold + centered
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108

# Set selection to match target
$ws.Range("E18").Select()